$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Formula = "=""66.774.58"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E2")
$c.Formula = "=""  +0.82%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D3")
$c.Formula = "=""3.236.82"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E3")
$c.Formula = "=""  +1.32%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D4")
$c.Formula = "=""1.00"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E4")
$c.Formula = "=""  +0.01%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D5")
$c.Formula = "=""606.24"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E5")
$c.Formula = "=""  +1.57%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D6")
$c.Formula = "=""157.54"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E6")
$c.Formula = "=""  +1.97%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D8")
$c.Formula = "=""3.237.17"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E8")
$c.Formula = "=""  +1.33%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D9")
$c.Formula = "=""0.548"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E9")
$c.Formula = "=""  +2.23%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E10")
$c.Formula = "=""  +0.46%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D11")
$c.Formula = "=""5.70"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E11")
$c.Formula = "=""  -6.64%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D12")
$c.Formula = "=""0.511"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E12")
$c.Formula = "=""  -0.66%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D13")
$c.Formula = "=""0.0000273"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E13")
$c.Formula = "=""  +0.89%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D14")
$c.Formula = "=""39.13"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E14")
$c.Formula = "=""  +0.01%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D15")
$c.Formula = "=""3.760.72"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E15")
$c.Formula = "=""  +1.14%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D16")
$c.Formula = "=""66.839.17"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E16")
$c.Formula = "=""  +1.03%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D17")
$c.Formula = "=""7.46"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E17")
$c.Formula = "=""  +0.28%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D18")
$c.Formula = "=""3.231.06"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E18")
$c.Formula = "=""  +0.99%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E19")
$c.Formula = "=""  +1.24%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D20")
$c.Formula = "=""513.13"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E20")
$c.Formula = "=""  +0.56%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D21")
$c.Formula = "=""15.42"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E21")
$c.Formula = "=""  +0.68%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D22")
$c.Formula = "=""0.740"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E22")
$c.Formula = "=""  +0.19%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D23")
$c.Formula = "=""8.12"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E23")
$c.Formula = "=""  +1.38%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D24")
$c.Formula = "=""14.90"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D25")
$c.Formula = "=""84.87"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E25")
$c.Formula = "=""  -0.01%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E26")
$c.Formula = "=""  +0.13%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D27")
$c.Formula = "=""9.56"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E27")
$c.Formula = "=""  +2.92%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D28")
$c.Formula = "=""3.02"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E28")
$c.Formula = "=""  +0.68%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D29")
$c.Formula = "=""2.42"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E29")
$c.Formula = "=""  +5.87%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D30")
$c.Formula = "=""3.07"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E30")
$c.Formula = "=""  +5.89%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D31")
$c.Formula = "=""7.11"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D32")
$c.Formula = "=""28.26"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E32")
$c.Formula = "=""  -0.14%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E33")
$c.Formula = "=""  +0.18%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E34")
$c.Formula = "=""  -3.25%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D35")
$c.Formula = "=""6.56"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E35")
$c.Formula = "=""  +0.01%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D36")
$c.Formula = "=""520.36"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E36")
$c.Formula = "=""  +7.37%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D37")
$c.Formula = "=""56.38"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E37")
$c.Formula = "=""  +2.75%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D38")
$c.Formula = "=""0.0927"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E38")
$c.Formula = "=""  +2.82%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D39")
$c.Formula = "=""0.0₃0765"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E39")
$c.Formula = "=""  +17.03%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D40")
$c.Formula = "=""0.0421"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E40")
$c.Formula = "=""  +0.66%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("B41")
$c.Formula = "=""dogwifhat"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("C41")
$c.Formula = "=""https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D41")
$c.Formula = "=""3.02"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E41")
$c.Formula = "=""  +3.26%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("B42")
$c.Formula = "=""Kaspa"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("C42")
$c.Formula = "=""https://coinranking.com/coin/V8GxkwWow+kaspa-kas"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D42")
$c.Formula = "=""0.128"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E42")
$c.Formula = "=""  +4.76%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D43")
$c.Formula = "=""8.83"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E43")
$c.Formula = "=""  -0.24%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D44")
$c.Formula = "=""0.304"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E44")
$c.Formula = "=""  +1.95%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D45")
$c.Formula = "=""2.54"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E45")
$c.Formula = "=""  +4.97%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D46")
$c.Formula = "=""2.862.33"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E46")
$c.Formula = "=""  -2.29%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D47")
$c.Formula = "=""28.59"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E47")
$c.Formula = "=""  +0.27%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E48")
$c.Formula = "=""  +4.07%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E49")
$c.Formula = "=""  -0.08%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E50")
$c.Formula = "=""  +0.48%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("D51")
$c.Formula = "=""2.62"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
$c = $ws.Range("E51")
$c.Formula = "=""  +1.10%  """
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null
